# 7.1_Input_PFHpA.xlsx - "data" sheet updates
#
# Commit: "May extra edits 2 / Extra edits code: K reduction fouling
# extensions, full-scale simulations"
#
# The underlying change is a refresh of the influent/effluent breakthrough
# data on the "data" worksheet: the time values for the two "174 day"
# samples move to "124 days", and the corresponding influent concentration
# values are corrected from 50000 to 0.04. The workbook is also left with
# the "data" sheet active/selected (cell C4), matching the state the
# workbook was saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Influent row at t=0: concentration 50000 -> 0.04
$ws.Range("C2").Value = 0.04

# Influent row at t=174 -> t=124, concentration 50000 -> 0.04
$ws.Range("B3").Value = 124
$ws.Range("C3").Value = 0.04

# Effluent row at t=174 -> t=124
$ws.Range("B5").Value = 124

# Leave the workbook with the "data" sheet active and C4 selected, which
# also makes "data" the saved active tab (and clears the previous
# tabSelected on "Kdata").
$ws.Range("C4").Select() | Out-Null
